# Update the "supervision" worksheet with the two supervision/degree
# entries that replace the previous (PhD / Professional Doctorate) rows.
# Row 2 -> Master's in Clinical Neuropsychology (Sara Silva Gómez, UVIU)
# Row 3 -> Master's in Psychology (Yenny Johanna Baron Londoño, U. El Bosque)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("supervision")

# Row 2
$ws.Range("A2").Value = "Máster en Neuropsicología Clínica"
$ws.Range("B2").Value = "2022-2023"
$ws.Range("C2").Value = "Sara Silva Gómez"
$ws.Range("D2").Value = "\href{https://www.universidadviu.com/co/}{Universidad Internacional de Valencia}, España"

# Row 3
$ws.Range("A3").Value = "Maestría en Psicología"
$ws.Range("B3").Value = "2019 - 2020"
$ws.Range("C3").Value = "Yenny Johanna Baron Londoño"
$ws.Range("D3").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia"

# Match the author's final selection position recorded in the sheet view.
$ws.Range("C3").Select()

$wb.Save()
